$wb = $excel.ActiveWorkbook

# --- Rename the first worksheet: "3ASY01_RNASeq" -> "Assay" ---
$wsAssay = $wb.Worksheets.Item(1)
$wsAssay.Name = "Assay"

# --- Rename the annotation table on the Assay sheet ---
# annotationTableSpicySloth85 -> annotationTableSpicySloth84
$loAssay = $wsAssay.ListObjects.Item(1)
$loAssay.Name = "annotationTableSpicySloth84"

# --- Update the literal "Table" metadata value on SwateTemplateMetadata (B6) to match ---
$wsMeta = $wb.Worksheets.Item("SwateTemplateMetadata")
$wsMeta.Range("B6").Value = "annotationTableSpicySloth84"

# --- Update selections / active sheet to match the saved view state ---
# Assay sheet keeps a new cell selection but is not the active tab.
$wsAssay.Range("X24").Select()

# SwateTemplateMetadata becomes the active tab with a new selection.
$wsMeta.Range("H33").Select()
